$d = $word.ActiveDocument

function Add-TrailingSpaceHyperlink([string]$urlText) {
    # Find the plain-text URL run (exact text, case-sensitive, whole content).
    $found = $d.Content.Duplicate
    $ok = $found.Find.Execute($urlText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $urlText"
        return
    }

    # Insert a literal space right after the URL text (this becomes its own run,
    # outside of the future hyperlink range).
    $insertPoint = $d.Range($found.End, $found.End)
    $insertPoint.InsertAfter(" ")

    # Re-establish a range over exactly the URL text (unaffected by the insert,
    # since the insert happened at/after $found.End) and turn it into a hyperlink.
    $linkRange = $d.Range($found.Start, $found.End)
    $d.Hyperlinks.Add($linkRange, $urlText) | Out-Null
}

# 1) "https://aka.ms/bcle" -> hyperlink (target https://aka.ms/bcle) + trailing space run
Add-TrailingSpaceHyperlink "https://aka.ms/bcle"

# 2) Blogs URL -> hyperlink + trailing space run
Add-TrailingSpaceHyperlink "https://github.com/microsoft/BCTech/blob/master/samples/AppInsights/BLOGS.md"

# 3) Videos URL -> hyperlink + trailing space run
Add-TrailingSpaceHyperlink "https://github.com/microsoft/BCTech/blob/master/samples/AppInsights/VIDEOS.md"

# 4) "April" -> "May" (the newsletter month reference near the end of the doc)
$d.Content.Find.Execute("April", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "May", 2) | Out-Null

Write-Output "done"
